$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.287.97"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.622.96"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.08"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.39"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.623.91"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.51"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.25"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.088.62"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.091.85"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000146"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.638.25"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.53"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.67"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.72"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.07"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.61"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "546.13"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.03"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0838"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.23"
$ws.Range("E36").Value = "  +2.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "168.93"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.402"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("E40").Value = "  +5.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.99"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.24"
$ws.Range("E43").Value = "  -5.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.80"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.76"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.623"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.95"
$ws.Range("E50").Value = "  +13.13%  "
$ws.Range("E51").Value = "  -0.71%  "
